$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55: A Real Smooth Move / Lanolin
$ws.Range("H55").Value = 473.5
$ws.Range("I55").Value = 474.66666
$ws.Range("J55").Value = 470
$ws.Range("K55").Value = 474.66666
$ws.Range("L55").Value = 470
$ws.Range("M55").Value = -260.66666
$ws.Range("N55").Value = -898

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 3492.2307
$ws.Range("I113").Value = 2979.8
$ws.Range("J113").Value = 3812.5
$ws.Range("K113").Value = 2979.8
$ws.Range("L113").Value = 3812.5
$ws.Range("M113").Value = 274.1999999999998
$ws.Range("N113").Value = -10320.5

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 2606956.5
$ws.Range("I116").Value = 2780456.8
$ws.Range("J116").Value = 4453
$ws.Range("K116").Value = 2780456.8
$ws.Range("L116").Value = 4453
$ws.Range("M116").Value = -2777014.8
$ws.Range("N116").Value = -11337

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4371.956
$ws.Range("I138").Value = 2450.3635
$ws.Range("J138").Value = 5290.978
$ws.Range("K138").Value = 7351.0905
$ws.Range("L138").Value = 15872.934
$ws.Range("M138").Value = -2211.0905
$ws.Range("N138").Value = -26152.934

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 3907.4285
$ws.Range("I141").Value = 2570.4
$ws.Range("J141").Value = 7250
$ws.Range("K141").Value = 7711.200000000001
$ws.Range("L141").Value = 21750
$ws.Range("M141").Value = -2531.200000000001
$ws.Range("N141").Value = -32110


$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 721.8823
$ws.Range("I2").Value = 642
$ws.Range("J2").Value = 981.5
$ws.Range("K2").Value = 642
$ws.Range("L2").Value = 981.5
$ws.Range("M2").Value = -529
$ws.Range("N2").Value = -1207.5

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2692.0625
$ws.Range("I74").Value = 2686.5925
$ws.Range("J74").Value = 2721.6
$ws.Range("K74").Value = 2686.5925
$ws.Range("L74").Value = 2721.6
$ws.Range("M74").Value = -1812.5925
$ws.Range("N74").Value = -4469.6

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2692.0625
$ws.Range("I77").Value = 2686.5925
$ws.Range("J77").Value = 2721.6
$ws.Range("K77").Value = 13432.9625
$ws.Range("L77").Value = 13608
$ws.Range("M77").Value = -9064.962500000001
$ws.Range("N77").Value = -22344

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 721.8823
$ws.Range("I116").Value = 642
$ws.Range("J116").Value = 981.5
$ws.Range("K116").Value = 642
$ws.Range("L116").Value = 981.5
$ws.Range("M116").Value = 1652
$ws.Range("N116").Value = -5569.5

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2676.2114
$ws.Range("I122").Value = 2440.561
$ws.Range("J122").Value = 3554.5454
$ws.Range("K122").Value = 7321.683000000001
$ws.Range("L122").Value = 10663.6362
$ws.Range("M122").Value = -4871.683000000001
$ws.Range("N122").Value = -15563.6362

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2505.5303
$ws.Range("I132").Value = 2316.6724
$ws.Range("K132").Value = 6950.0172
$ws.Range("M132").Value = -4420.0172


$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 721.8823
$ws.Range("I3").Value = 642
$ws.Range("J3").Value = 981.5
$ws.Range("K3").Value = 642
$ws.Range("L3").Value = 981.5
$ws.Range("M3").Value = -528
$ws.Range("N3").Value = -1209.5

# Row 59: Pop That Top / Cobalt Raising Hammer
$ws.Range("H59").Value = 20390
$ws.Range("J59").Value = 20390
$ws.Range("L59").Value = 20390
$ws.Range("N59").Value = -22084

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1619.091
$ws.Range("I107").Value = 1356.3334
$ws.Range("J107").Value = 2319.7778
$ws.Range("K107").Value = 1356.3334
$ws.Range("L107").Value = 2319.7778
$ws.Range("M107").Value = 563.6666
$ws.Range("N107").Value = -6159.7778

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1398.5862
$ws.Range("I134").Value = 1161.091
$ws.Range("J134").Value = 2145
$ws.Range("K134").Value = 3483.273
$ws.Range("L134").Value = 6435
$ws.Range("M134").Value = -948.2729999999997
$ws.Range("N134").Value = -11505


$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1815.3939
$ws.Range("I16").Value = 1674.3704
$ws.Range("J16").Value = 2450
$ws.Range("K16").Value = 1674.3704
$ws.Range("L16").Value = 2450
$ws.Range("M16").Value = -1387.3704
$ws.Range("N16").Value = -3024

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2727.389
$ws.Range("I31").Value = 1407.3334
$ws.Range("J31").Value = 4047.4443
$ws.Range("K31").Value = 1407.3334
$ws.Range("L31").Value = 4047.4443
$ws.Range("M31").Value = -1112.3334
$ws.Range("N31").Value = -4637.4443

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2727.389
$ws.Range("I34").Value = 1407.3334
$ws.Range("J34").Value = 4047.4443
$ws.Range("K34").Value = 1407.3334
$ws.Range("L34").Value = 4047.4443
$ws.Range("M34").Value = -1205.3334
$ws.Range("N34").Value = -4451.4443

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2745.9119
$ws.Range("I58").Value = 894.913
$ws.Range("J58").Value = 6616.1816
$ws.Range("K58").Value = 894.913
$ws.Range("L58").Value = 6616.1816
$ws.Range("M58").Value = -691.913
$ws.Range("N58").Value = -7022.1816

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 3649.5
$ws.Range("I99").Value = 3826
$ws.Range("J99").Value = 3120
$ws.Range("K99").Value = 3826
$ws.Range("L99").Value = 3120
$ws.Range("M99").Value = -2328
$ws.Range("N99").Value = -6116

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1016.3333
$ws.Range("I107").Value = 856.64703
$ws.Range("J107").Value = 1695
$ws.Range("K107").Value = 856.64703
$ws.Range("L107").Value = 1695
$ws.Range("M107").Value = 1063.35297
$ws.Range("N107").Value = -5535

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1815.3939
$ws.Range("I113").Value = 1674.3704
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 1674.3704
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = 495.6296
$ws.Range("N113").Value = -6790

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 957.63635
$ws.Range("I122").Value = 915.3333
$ws.Range("J122").Value = 1008.4
$ws.Range("K122").Value = 2745.9999
$ws.Range("L122").Value = 3025.2
$ws.Range("M122").Value = -295.9998999999998
$ws.Range("N122").Value = -7925.2

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 3649.5
$ws.Range("I126").Value = 3826
$ws.Range("J126").Value = 3120
$ws.Range("K126").Value = 11478
$ws.Range("L126").Value = 9360
$ws.Range("M126").Value = -9008
$ws.Range("N126").Value = -14300

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1429.7576
$ws.Range("I132").Value = 1145.5385
$ws.Range("J132").Value = 2485.4285
$ws.Range("K132").Value = 3436.6155
$ws.Range("L132").Value = 7456.2855
$ws.Range("M132").Value = -906.6155000000003
$ws.Range("N132").Value = -12516.2855

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2012.0754
$ws.Range("I134").Value = 2163.348
$ws.Range("K134").Value = 6490.044
$ws.Range("M134").Value = -3955.044

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2745.9119
$ws.Range("I136").Value = 894.913
$ws.Range("J136").Value = 6616.1816
$ws.Range("K136").Value = 2684.739
$ws.Range("L136").Value = 19848.5448
$ws.Range("M136").Value = -134.739
$ws.Range("N136").Value = -24948.5448


$ws = $wb.Worksheets.Item("CUL")
# Row 50: Moving Up in the World / Rolanberry Cheese
$ws.Range("H50").Value = 117.916664
$ws.Range("I50").Value = 31.6
$ws.Range("J50").Value = 549.5
$ws.Range("K50").Value = 94.80000000000001
$ws.Range("L50").Value = 1648.5
$ws.Range("M50").Value = 386.2
$ws.Range("N50").Value = -2610.5

# Row 53: Rolanberry Fields Forever / Rolanberry Cheese
$ws.Range("H53").Value = 117.916664
$ws.Range("I53").Value = 31.6
$ws.Range("J53").Value = 549.5
$ws.Range("K53").Value = 94.80000000000001
$ws.Range("L53").Value = 1648.5
$ws.Range("M53").Value = 386.2
$ws.Range("N53").Value = -2610.5

# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 563.6667
$ws.Range("I92").Value = 495.5
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 1486.5
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -238.5

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 873.09375
$ws.Range("I131").Value = 492.92307
$ws.Range("J131").Value = 1133.2106
$ws.Range("K131").Value = 1478.76921
$ws.Range("L131").Value = 3399.6318
$ws.Range("M131").Value = 3561.23079
$ws.Range("N131").Value = -13479.6318

# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 1443.871
$ws.Range("I140").Value = 1487.3684
$ws.Range("J140").Value = 1375
$ws.Range("K140").Value = 4462.1052
$ws.Range("L140").Value = 4125
$ws.Range("M140").Value = 717.8948
$ws.Range("N140").Value = -14485


$ws = $wb.Worksheets.Item("GSM")
# Row 93: One Ring Circus / Triphane Ring of Slaying
$ws.Range("H93").Value = 6666.6665
$ws.Range("J93").Value = 6666.6665
$ws.Range("L93").Value = 6666.6665
$ws.Range("N93").Value = -10410.6665


$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 310.89655
$ws.Range("I55").Value = 368.13333
$ws.Range("J55").Value = 249.57143
$ws.Range("K55").Value = 368.13333
$ws.Range("L55").Value = 249.57143
$ws.Range("M55").Value = -195.13333
$ws.Range("N55").Value = -595.57143

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3699.125
$ws.Range("I122").Value = 4024.65
$ws.Range("J122").Value = 3156.5833
$ws.Range("K122").Value = 12073.95
$ws.Range("L122").Value = 9469.749899999999
$ws.Range("M122").Value = -9623.950000000001
$ws.Range("N122").Value = -14369.7499

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 7956.1816
$ws.Range("I136").Value = 8813.5625
$ws.Range("K136").Value = 26440.6875
$ws.Range("M136").Value = -23890.6875


$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 858.46155
$ws.Range("I122").Value = 777.7778
$ws.Range("J122").Value = 1040
$ws.Range("K122").Value = 2333.3334
$ws.Range("L122").Value = 3120
$ws.Range("M122").Value = 116.6666
$ws.Range("N122").Value = -8020

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1811.8667
$ws.Range("I132").Value = 923.87177
$ws.Range("J132").Value = 3461
$ws.Range("K132").Value = 2771.61531
$ws.Range("L132").Value = 10383
$ws.Range("M132").Value = -241.6153100000001
$ws.Range("N132").Value = -15443

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 6642.8096
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 16428.428
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 49285.284
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -54385.284


# CUL row 92: add new N92 cell (LeveProfitHQ) that did not exist before
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N92").Value = -4596
